$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "ICP AA / ICP AG / ICP NAV" values from N24:P24.
# (Once unreferenced, the backing shared-string entries are dropped and
#  every later shared-string index renumbers down by 3 automatically -
#  matching the rest of the sheet's v="..." shifts in the diff.)
$ws.Cells.Item(24, 14).Value = ""
$ws.Cells.Item(24, 15).Value = ""
$ws.Cells.Item(24, 16).Value = ""

# New "ICP Buttons" mini-table added in columns Q:R, rows 2-12.
$ws.Range("Q2").Value = "ICP Buttons"
$ws.Range("R2").Value = "F1-F7"

$ws.Range("Q3").Value = "DRIFT C/O"
$ws.Range("R3").Value = "F8"

$ws.Range("Q4").Value = "DRIFT C/O WARN RESET"
$ws.Range("R4").Value = "F9"

$ws.Range("Q5").Value = "FLCS RESET"
$ws.Range("R5").Value = "F10"

$ws.Range("Q6").Value = "Steerpoints"
$ws.Range("R6").Value = "<, >"

$ws.Range("Q7").Value = "HSI HDG Knob"
$ws.Range("R7").Value = "Ins, Del"

# R8 is written before Q8 so the shared-string table picks up
# "Home, End" ahead of "HSI CRS Knob" (matches the source order).
$ws.Range("R8").Value = "Home, End"
$ws.Range("Q8").Value = "HSI CRS Knob"

$ws.Range("Q9").Value = "Pressure Knob"
$ws.Range("R9").Value = "PgUp, PgDn"

$ws.Range("Q10").Value = "INSTR MODE Cycle"
$ws.Range("R10").Value = "\"

$ws.Range("Q11").Value = "CMDS PRGM Knob"
$ws.Range("R11").Value = "[ ,]"

$ws.Range("Q12").Value = "CMDS MODE Knob"
$ws.Range("R12").Value = ';, '''

# Widen column Q to fit the new labels.
$ws.Columns.Item(17).ColumnWidth = 21.7109375

# Move the active selection to N19.
[void]$ws.Range("N19").Select()
